$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct two pre-existing values in row 199
$ws.Cells.Item(199, 7).Value2 = 0.5446734   # G199
$ws.Cells.Item(199, 21).Value2 = 0.6103289  # U199

# Row 200 (17 08 2020)
$ws.Cells.Item(200, 1).Value2 = "17 08 2020"
$ws.Cells.Item(200, 2).Value2 = 0.4047853
$ws.Cells.Item(200, 3).Value2 = 1.1029643
$ws.Cells.Item(200, 4).Value2 = 0.6807266
$ws.Cells.Item(200, 6).Value2 = 0.7336025
$ws.Cells.Item(200, 7).Value2 = 0.5477414
$ws.Cells.Item(200, 8).Value2 = 0.4441322
$ws.Cells.Item(200, 9).Value2 = 0.1490105
$ws.Cells.Item(200, 10).Value2 = 0.2616747
$ws.Cells.Item(200, 11).Value2 = 0.3301592
$ws.Cells.Item(200, 12).Value2 = 0.7173486
$ws.Cells.Item(200, 13).Value2 = 0.9665268
$ws.Cells.Item(200, 15).Value2 = 0.5133129
$ws.Cells.Item(200, 16).Value2 = 0.6639978
$ws.Cells.Item(200, 17).Value2 = 0.7683105
$ws.Cells.Item(200, 18).Value2 = 0.4261091
$ws.Cells.Item(200, 19).Value2 = 0.798144
$ws.Cells.Item(200, 20).Value2 = 0.6231455
$ws.Cells.Item(200, 21).Value2 = 0.6128055
$ws.Cells.Item(200, 22).Value2 = 1.0070671
$ws.Cells.Item(200, 23).Value2 = 0.2660783
$ws.Cells.Item(200, 24).Value2 = 0.3578887
$ws.Cells.Item(200, 25).Value2 = 0.2237744
$ws.Cells.Item(200, 26).Value2 = 0.4165883
$ws.Cells.Item(200, 27).Value2 = 0.4319993
$ws.Cells.Item(200, 28).Value2 = 0.5782555
$ws.Cells.Item(200, 30).Value2 = 1.1105631
$ws.Cells.Item(200, 31).Value2 = 0.5049285
$ws.Cells.Item(200, 32).Value2 = 0.453165
$ws.Cells.Item(200, 33).Value2 = 0.6911835
$ws.Cells.Item(200, 34).Value2 = 0.680644
$ws.Cells.Item(200, 35).Value2 = 0.1889782
$ws.Cells.Item(200, 36).Value2 = 0.2660204
$ws.Cells.Item(200, 37).Value2 = 0.4510177
$ws.Cells.Item(200, 38).Value2 = 0.7194225
$ws.Cells.Item(200, 39).Value2 = 0.3143743
$ws.Cells.Item(200, 40).Value2 = 0.5443116
$ws.Cells.Item(200, 41).Value2 = 0.6590832
$ws.Cells.Item(200, 42).Value2 = 0.3919441
$ws.Cells.Item(200, 43).Value2 = 0.3943336
$ws.Cells.Item(200, 45).Value2 = 0.3545821
$ws.Cells.Item(200, 46).Value2 = 0.7990294999999999
$ws.Cells.Item(200, 47).Value2 = 0.3915695
$ws.Cells.Item(200, 48).Value2 = 0.6717266
$ws.Cells.Item(200, 49).Value2 = 0.6926722
$ws.Cells.Item(200, 50).Value2 = 0.692263
$ws.Cells.Item(200, 51).Value2 = 0.4782376
$ws.Cells.Item(200, 53).Value2 = 0.1392754
$ws.Cells.Item(200, 54).Value2 = 0.4137513
$ws.Cells.Item(200, 55).Value2 = 0.4275651
$ws.Cells.Item(200, 56).Value2 = 0.3814741
$ws.Cells.Item(200, 57).Value2 = 0.3808559

# Row 201 (18 08 2020)
$ws.Cells.Item(201, 1).Value2 = "18 08 2020"
$ws.Cells.Item(201, 2).Value2 = 0.4921588
$ws.Cells.Item(201, 3).Value2 = 1.0643424
$ws.Cells.Item(201, 4).Value2 = 0.6356991
$ws.Cells.Item(201, 6).Value2 = 0.6563658999999999
$ws.Cells.Item(201, 7).Value2 = 0.5165051000000001
$ws.Cells.Item(201, 8).Value2 = 0.4238525
$ws.Cells.Item(201, 9).Value2 = 0.1445852
$ws.Cells.Item(201, 10).Value2 = 0.1230012
$ws.Cells.Item(201, 11).Value2 = 0.3613161
$ws.Cells.Item(201, 12).Value2 = 0.686524
$ws.Cells.Item(201, 13).Value2 = 1.0113092
$ws.Cells.Item(201, 15).Value2 = 0.4808618
$ws.Cells.Item(201, 16).Value2 = 0.7981363
$ws.Cells.Item(201, 17).Value2 = 0.7951613
$ws.Cells.Item(201, 18).Value2 = 0.3851701
$ws.Cells.Item(201, 19).Value2 = 0.7960521
$ws.Cells.Item(201, 20).Value2 = 0.7347946
$ws.Cells.Item(201, 21).Value2 = 0.6360099
$ws.Cells.Item(201, 22).Value2 = 1.0192065
$ws.Cells.Item(201, 23).Value2 = 0.2693894
$ws.Cells.Item(201, 24).Value2 = 0.383336
$ws.Cells.Item(201, 25).Value2 = 0.213134
$ws.Cells.Item(201, 26).Value2 = 0.4490157
$ws.Cells.Item(201, 27).Value2 = 0.3885677
$ws.Cells.Item(201, 28).Value2 = 0.62588
$ws.Cells.Item(201, 30).Value2 = 1.0763115
$ws.Cells.Item(201, 31).Value2 = 0.5664476000000001
$ws.Cells.Item(201, 32).Value2 = 0.4807025
$ws.Cells.Item(201, 33).Value2 = 0.6790136
$ws.Cells.Item(201, 34).Value2 = 0.6508283
$ws.Cells.Item(201, 35).Value2 = 0.2182613
$ws.Cells.Item(201, 36).Value2 = 0.2332097
$ws.Cells.Item(201, 37).Value2 = 0.3861522
$ws.Cells.Item(201, 38).Value2 = 0.6263351
$ws.Cells.Item(201, 39).Value2 = 0.3118155
$ws.Cells.Item(201, 40).Value2 = 0.5439827
$ws.Cells.Item(201, 41).Value2 = 0.6369738
$ws.Cells.Item(201, 42).Value2 = 0.4281315
$ws.Cells.Item(201, 43).Value2 = 0.390042
$ws.Cells.Item(201, 45).Value2 = 0.3670048
$ws.Cells.Item(201, 46).Value2 = 0.7589256
$ws.Cells.Item(201, 47).Value2 = 0.5472068
$ws.Cells.Item(201, 48).Value2 = 0.6903181
$ws.Cells.Item(201, 49).Value2 = 0.6837459
$ws.Cells.Item(201, 50).Value2 = 0.6407502
$ws.Cells.Item(201, 51).Value2 = 0.446654
$ws.Cells.Item(201, 53).Value2 = 0.1757628
$ws.Cells.Item(201, 54).Value2 = 0.4068253
$ws.Cells.Item(201, 55).Value2 = 0.4831073
$ws.Cells.Item(201, 56).Value2 = 0.4451979
$ws.Cells.Item(201, 57).Value2 = 0.50065

# Row 202 (19 08 2020)
$ws.Cells.Item(202, 1).Value2 = "19 08 2020"
$ws.Cells.Item(202, 2).Value2 = 0.555
$ws.Cells.Item(202, 3).Value2 = 1.0669742
$ws.Cells.Item(202, 4).Value2 = 0.5781856
$ws.Cells.Item(202, 6).Value2 = 0.6653272
$ws.Cells.Item(202, 7).Value2 = 0.5322966
$ws.Cells.Item(202, 8).Value2 = 0.410507
$ws.Cells.Item(202, 9).Value2 = 0.1779487
$ws.Cells.Item(202, 10).Value2 = 0.1240695
$ws.Cells.Item(202, 11).Value2 = 0.306788
$ws.Cells.Item(202, 12).Value2 = 0.7040899
$ws.Cells.Item(202, 13).Value2 = 0.8874875
$ws.Cells.Item(202, 15).Value2 = 0.4887067
$ws.Cells.Item(202, 16).Value2 = 0.8437114999999999
$ws.Cells.Item(202, 17).Value2 = 0.9578727
$ws.Cells.Item(202, 18).Value2 = 0.4132225
$ws.Cells.Item(202, 19).Value2 = 0.7855347
$ws.Cells.Item(202, 20).Value2 = 0.7426700000000001
$ws.Cells.Item(202, 21).Value2 = 0.6882355999999999
$ws.Cells.Item(202, 22).Value2 = 1.1224403
$ws.Cells.Item(202, 23).Value2 = 0.2774423
$ws.Cells.Item(202, 24).Value2 = 0.3480744
$ws.Cells.Item(202, 25).Value2 = 0.2240533
$ws.Cells.Item(202, 26).Value2 = 0.4703502
$ws.Cells.Item(202, 27).Value2 = 0.3608298
$ws.Cells.Item(202, 28).Value2 = 0.6657531
$ws.Cells.Item(202, 30).Value2 = 1.0948521
$ws.Cells.Item(202, 31).Value2 = 0.5579843
$ws.Cells.Item(202, 32).Value2 = 0.4968115
$ws.Cells.Item(202, 33).Value2 = 0.7347013999999999
$ws.Cells.Item(202, 34).Value2 = 0.7234692
$ws.Cells.Item(202, 35).Value2 = 0.3012102
$ws.Cells.Item(202, 36).Value2 = 0.2556868
$ws.Cells.Item(202, 37).Value2 = 0.4039599
$ws.Cells.Item(202, 38).Value2 = 0.5959876
$ws.Cells.Item(202, 39).Value2 = 0.3285346
$ws.Cells.Item(202, 40).Value2 = 0.4767433
$ws.Cells.Item(202, 41).Value2 = 0.6641992
$ws.Cells.Item(202, 42).Value2 = 0.4525957
$ws.Cells.Item(202, 43).Value2 = 0.3477199
$ws.Cells.Item(202, 45).Value2 = 0.4329054
$ws.Cells.Item(202, 46).Value2 = 0.7585933
$ws.Cells.Item(202, 47).Value2 = 0.4441092
$ws.Cells.Item(202, 48).Value2 = 0.6927337
$ws.Cells.Item(202, 49).Value2 = 0.6852665999999999
$ws.Cells.Item(202, 50).Value2 = 0.6927569
$ws.Cells.Item(202, 51).Value2 = 0.4587121
$ws.Cells.Item(202, 53).Value2 = 0.08498559999999999
$ws.Cells.Item(202, 54).Value2 = 0.3723756
$ws.Cells.Item(202, 55).Value2 = 0.4190395
$ws.Cells.Item(202, 56).Value2 = 0.5484771000000001
$ws.Cells.Item(202, 57).Value2 = 0.4162067

# Row 203 (20 08 2020)
$ws.Cells.Item(203, 1).Value2 = "20 08 2020"
$ws.Cells.Item(203, 2).Value2 = 0.7377939999999999
$ws.Cells.Item(203, 3).Value2 = 1.0597842
$ws.Cells.Item(203, 4).Value2 = 0.7647505
$ws.Cells.Item(203, 6).Value2 = 0.7106584
$ws.Cells.Item(203, 7).Value2 = 0.5266883
$ws.Cells.Item(203, 8).Value2 = 0.4227711
$ws.Cells.Item(203, 9).Value2 = 0.1705545
$ws.Cells.Item(203, 10).Value2 = 0.1278772
$ws.Cells.Item(203, 11).Value2 = 0.3092179
$ws.Cells.Item(203, 12).Value2 = 0.6498813
$ws.Cells.Item(203, 13).Value2 = 0.9110343
$ws.Cells.Item(203, 15).Value2 = 0.4356543
$ws.Cells.Item(203, 16).Value2 = 0.88866
$ws.Cells.Item(203, 17).Value2 = 1.0728105
$ws.Cells.Item(203, 18).Value2 = 0.4238124
$ws.Cells.Item(203, 19).Value2 = 0.7387834
$ws.Cells.Item(203, 20).Value2 = 0.7949887
$ws.Cells.Item(203, 21).Value2 = 0.6673875
$ws.Cells.Item(203, 22).Value2 = 1.1693106
$ws.Cells.Item(203, 23).Value2 = 0.308831
$ws.Cells.Item(203, 24).Value2 = 0.3614437
$ws.Cells.Item(203, 25).Value2 = 0.211205
$ws.Cells.Item(203, 26).Value2 = 0.4581249
$ws.Cells.Item(203, 27).Value2 = 0.3574208
$ws.Cells.Item(203, 28).Value2 = 0.6351087
$ws.Cells.Item(203, 30).Value2 = 1.0972169
$ws.Cells.Item(203, 31).Value2 = 0.6391864
$ws.Cells.Item(203, 32).Value2 = 0.5167403
$ws.Cells.Item(203, 33).Value2 = 0.6422848
$ws.Cells.Item(203, 34).Value2 = 0.8304295
$ws.Cells.Item(203, 35).Value2 = 0.3285427
$ws.Cells.Item(203, 36).Value2 = 0.2759627
$ws.Cells.Item(203, 37).Value2 = 0.4631256
$ws.Cells.Item(203, 38).Value2 = 0.5479685
$ws.Cells.Item(203, 39).Value2 = 0.3210486
$ws.Cells.Item(203, 40).Value2 = 0.4532837
$ws.Cells.Item(203, 41).Value2 = 0.719542
$ws.Cells.Item(203, 42).Value2 = 0.3929036
$ws.Cells.Item(203, 43).Value2 = 0.3352904
$ws.Cells.Item(203, 45).Value2 = 0.372028
$ws.Cells.Item(203, 46).Value2 = 0.7966444
$ws.Cells.Item(203, 47).Value2 = 0.7538283
$ws.Cells.Item(203, 48).Value2 = 0.6397855
$ws.Cells.Item(203, 49).Value2 = 0.7054222
$ws.Cells.Item(203, 50).Value2 = 0.6306071
$ws.Cells.Item(203, 51).Value2 = 0.4779283
$ws.Cells.Item(203, 53).Value2 = 0.1572923
$ws.Cells.Item(203, 54).Value2 = 0.3579087
$ws.Cells.Item(203, 55).Value2 = 0.4056082
$ws.Cells.Item(203, 56).Value2 = 0.5990788
$ws.Cells.Item(203, 57).Value2 = 0.6501186

# Row 204 (21 08 2020)
$ws.Cells.Item(204, 1).Value2 = "21 08 2020"
$ws.Cells.Item(204, 2).Value2 = 0.65123
$ws.Cells.Item(204, 3).Value2 = 1.0699572
$ws.Cells.Item(204, 4).Value2 = 0.7288888
$ws.Cells.Item(204, 6).Value2 = 0.630291
$ws.Cells.Item(204, 7).Value2 = 0.5113987
$ws.Cells.Item(204, 8).Value2 = 0.4856735
$ws.Cells.Item(204, 9).Value2 = 0.2345475
$ws.Cells.Item(204, 10).Value2 = 0.1329787
$ws.Cells.Item(204, 11).Value2 = 0.2729636
$ws.Cells.Item(204, 12).Value2 = 0.6545481
$ws.Cells.Item(204, 13).Value2 = 0.928884
$ws.Cells.Item(204, 15).Value2 = 0.4718965
$ws.Cells.Item(204, 16).Value2 = 0.9477365
$ws.Cells.Item(204, 17).Value2 = 1.1659769
$ws.Cells.Item(204, 18).Value2 = 0.461663
$ws.Cells.Item(204, 19).Value2 = 0.781073
$ws.Cells.Item(204, 20).Value2 = 0.773631
$ws.Cells.Item(204, 21).Value2 = 0.6864748000000001
$ws.Cells.Item(204, 22).Value2 = 1.1583645
$ws.Cells.Item(204, 23).Value2 = 0.2708283
$ws.Cells.Item(204, 24).Value2 = 0.3920199
$ws.Cells.Item(204, 25).Value2 = 0.2142585
$ws.Cells.Item(204, 26).Value2 = 0.5057717
$ws.Cells.Item(204, 27).Value2 = 0.3088636
$ws.Cells.Item(204, 28).Value2 = 0.5964187
$ws.Cells.Item(204, 30).Value2 = 1.1046447
$ws.Cells.Item(204, 31).Value2 = 0.6742527
$ws.Cells.Item(204, 32).Value2 = 0.5175343
$ws.Cells.Item(204, 33).Value2 = 0.6666196
$ws.Cells.Item(204, 34).Value2 = 0.979167
$ws.Cells.Item(204, 35).Value2 = 0.4408945
$ws.Cells.Item(204, 36).Value2 = 0.2564421
$ws.Cells.Item(204, 37).Value2 = 0.43656
$ws.Cells.Item(204, 38).Value2 = 0.4816024
$ws.Cells.Item(204, 39).Value2 = 0.3231496
$ws.Cells.Item(204, 40).Value2 = 0.4517217
$ws.Cells.Item(204, 41).Value2 = 0.7899414
$ws.Cells.Item(204, 42).Value2 = 0.4248582
$ws.Cells.Item(204, 43).Value2 = 0.3389427
$ws.Cells.Item(204, 45).Value2 = 0.2842144
$ws.Cells.Item(204, 46).Value2 = 0.7385401
$ws.Cells.Item(204, 47).Value2 = 0.9450001
$ws.Cells.Item(204, 48).Value2 = 0.6952983
$ws.Cells.Item(204, 49).Value2 = 0.7322689
$ws.Cells.Item(204, 50).Value2 = 0.6006601
$ws.Cells.Item(204, 51).Value2 = 0.5079597
$ws.Cells.Item(204, 53).Value2 = 0.1089014
$ws.Cells.Item(204, 54).Value2 = 0.3023158
$ws.Cells.Item(204, 55).Value2 = 0.4445865
$ws.Cells.Item(204, 56).Value2 = 0.6322479
$ws.Cells.Item(204, 57).Value2 = 0.7814687

# Row 205 (22 08 2020)
$ws.Cells.Item(205, 1).Value2 = "22 08 2020"
$ws.Cells.Item(205, 2).Value2 = 0.8014152
$ws.Cells.Item(205, 3).Value2 = 1.06837
$ws.Cells.Item(205, 4).Value2 = 0.7922736
$ws.Cells.Item(205, 6).Value2 = 0.6573366
$ws.Cells.Item(205, 7).Value2 = 0.5018521
$ws.Cells.Item(205, 8).Value2 = 0.4956143
$ws.Cells.Item(205, 9).Value2 = 0.2485632
$ws.Cells.Item(205, 10).Value2 = 0.1338688
$ws.Cells.Item(205, 11).Value2 = 0.2725721
$ws.Cells.Item(205, 12).Value2 = 0.624416
$ws.Cells.Item(205, 13).Value2 = 0.8870452
$ws.Cells.Item(205, 15).Value2 = 0.2643852
$ws.Cells.Item(205, 16).Value2 = 0.9946002
$ws.Cells.Item(205, 17).Value2 = 1.1663578
$ws.Cells.Item(205, 18).Value2 = 0.4896874
$ws.Cells.Item(205, 19).Value2 = 0.7848851999999999
$ws.Cells.Item(205, 20).Value2 = 0.7872187
$ws.Cells.Item(205, 21).Value2 = 0.696858
$ws.Cells.Item(205, 22).Value2 = 1.0141159
$ws.Cells.Item(205, 23).Value2 = 0.2301138
$ws.Cells.Item(205, 24).Value2 = 0.3927205
$ws.Cells.Item(205, 25).Value2 = 0.3077399
$ws.Cells.Item(205, 26).Value2 = 0.4518918
$ws.Cells.Item(205, 27).Value2 = 0.3384049
$ws.Cells.Item(205, 28).Value2 = 0.6042568
$ws.Cells.Item(205, 30).Value2 = 1.0064998
$ws.Cells.Item(205, 31).Value2 = 0.7025865
$ws.Cells.Item(205, 32).Value2 = 0.550458
$ws.Cells.Item(205, 33).Value2 = 0.5513764
$ws.Cells.Item(205, 34).Value2 = 0.9974608
$ws.Cells.Item(205, 35).Value2 = 0.3771443
$ws.Cells.Item(205, 36).Value2 = 0.2706995
$ws.Cells.Item(205, 37).Value2 = 0.4189135
$ws.Cells.Item(205, 38).Value2 = 0.5354934
$ws.Cells.Item(205, 39).Value2 = 0.3691699
$ws.Cells.Item(205, 40).Value2 = 0.4506075
$ws.Cells.Item(205, 41).Value2 = 0.7760754
$ws.Cells.Item(205, 42).Value2 = 0.4722772
$ws.Cells.Item(205, 43).Value2 = 0.3842036
$ws.Cells.Item(205, 45).Value2 = 0.2913385
$ws.Cells.Item(205, 46).Value2 = 0.7355263
$ws.Cells.Item(205, 47).Value2 = 1.0727605
$ws.Cells.Item(205, 48).Value2 = 0.7000754
$ws.Cells.Item(205, 49).Value2 = 0.7062657
$ws.Cells.Item(205, 50).Value2 = 0.5482112
$ws.Cells.Item(205, 51).Value2 = 0.4335884
$ws.Cells.Item(205, 53).Value2 = 0.2511845
$ws.Cells.Item(205, 54).Value2 = 0.262574
$ws.Cells.Item(205, 55).Value2 = 0.5189091
$ws.Cells.Item(205, 56).Value2 = 0.612158
$ws.Cells.Item(205, 57).Value2 = 0.6845343

# Row 206 (23 08 2020)
$ws.Cells.Item(206, 1).Value2 = "23 08 2020"
$ws.Cells.Item(206, 2).Value2 = 0.8908825
$ws.Cells.Item(206, 3).Value2 = 1.0736486
$ws.Cells.Item(206, 4).Value2 = 0.7679851
$ws.Cells.Item(206, 6).Value2 = 0.6337133
$ws.Cells.Item(206, 7).Value2 = 0.4930303
$ws.Cells.Item(206, 8).Value2 = 0.4991314
$ws.Cells.Item(206, 9).Value2 = 0.2952216
$ws.Cells.Item(206, 10).Value2 = 0
$ws.Cells.Item(206, 11).Value2 = 0.2861759
$ws.Cells.Item(206, 12).Value2 = 0.6193899
$ws.Cells.Item(206, 13).Value2 = 0.8989995
$ws.Cells.Item(206, 15).Value2 = 0.3507653
$ws.Cells.Item(206, 16).Value2 = 0.8887589
$ws.Cells.Item(206, 17).Value2 = 1.1208659
$ws.Cells.Item(206, 18).Value2 = 0.5257509
$ws.Cells.Item(206, 19).Value2 = 0.6990634999999999
$ws.Cells.Item(206, 20).Value2 = 0.8043928
$ws.Cells.Item(206, 21).Value2 = 0.7371695
$ws.Cells.Item(206, 22).Value2 = 1.0766428
$ws.Cells.Item(206, 23).Value2 = 0.2352608
$ws.Cells.Item(206, 24).Value2 = 0.4088317
$ws.Cells.Item(206, 25).Value2 = 0.2881312
$ws.Cells.Item(206, 26).Value2 = 0.4685418
$ws.Cells.Item(206, 27).Value2 = 0.3617011
$ws.Cells.Item(206, 28).Value2 = 0.5403767
$ws.Cells.Item(206, 30).Value2 = 0.9120764
$ws.Cells.Item(206, 31).Value2 = 0.6953885
$ws.Cells.Item(206, 32).Value2 = 0.5078658
$ws.Cells.Item(206, 33).Value2 = 0.566796
$ws.Cells.Item(206, 34).Value2 = 0.9674160000000001
$ws.Cells.Item(206, 35).Value2 = 0.4061767
$ws.Cells.Item(206, 36).Value2 = 0.292115
$ws.Cells.Item(206, 37).Value2 = 0.3989687
$ws.Cells.Item(206, 38).Value2 = 0.6390883000000001
$ws.Cells.Item(206, 39).Value2 = 0.3677877
$ws.Cells.Item(206, 40).Value2 = 0.4552416
$ws.Cells.Item(206, 41).Value2 = 0.794315
$ws.Cells.Item(206, 42).Value2 = 0.4938191
$ws.Cells.Item(206, 43).Value2 = 0.3593729
$ws.Cells.Item(206, 45).Value2 = 0.1654133
$ws.Cells.Item(206, 46).Value2 = 0.7368405
$ws.Cells.Item(206, 47).Value2 = 1.1797567
$ws.Cells.Item(206, 48).Value2 = 0.7682007
$ws.Cells.Item(206, 49).Value2 = 0.6966457
$ws.Cells.Item(206, 50).Value2 = 0.5066854
$ws.Cells.Item(206, 51).Value2 = 0.4864584
$ws.Cells.Item(206, 53).Value2 = 0.3496164
$ws.Cells.Item(206, 54).Value2 = 0.2605408
$ws.Cells.Item(206, 55).Value2 = 0.5019535000000001
$ws.Cells.Item(206, 56).Value2 = 0.5710781
$ws.Cells.Item(206, 57).Value2 = 0.8382695

# Row 207 (24 08 2020)
$ws.Cells.Item(207, 1).Value2 = "24 08 2020"

# Row 208 (25 08 2020)
$ws.Cells.Item(208, 1).Value2 = "25 08 2020"

